$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BoM")
$dnf = $wb.Worksheets.Item("DNF")

# --- BoM sheet: update "Net Name" / "Net Label" columns (AB/AC) ---

# Row 10 - C2 (Unpolarized capacitor, Earth<->U1-UCAP net)
$bom.Range("AB10").Value = "Net-(U1-UCAP),Earth"
$bom.Range("AC10").Value = "Net-(U1-UCAP),Earth"

# Row 15 - now has two different values for Net Name vs Net Label
$bom.Range("AB15").Value = "/RESET2,/SCK2,+5V,/MISO2,/MOSI2,GND"
$bom.Range("AC15").Value = "MOSI2,GND"

# Row 16 - J3 connector pin net reorder
$bom.Range("AB16").Value = "Net-(J3-Pin_2),Net-(J3-Pin_5),Net-(J3-Pin_1),Net-(J3-Pin_4),Net-(J3-Pin_3)"
$bom.Range("AC16").Value = "Net-(J3-Pin_2),Net-(J3-Pin_5),Net-(J3-Pin_1),Net-(J3-Pin_4),Net-(J3-Pin_3)"

# Row 17 - J6 connector pin net reorder
$bom.Range("AB17").Value = "Net-(J6-Pin_1),Net-(J6-Pin_6),Net-(J6-Pin_2),Net-(J6-Pin_5),Net-(J6-Pin_3),Net-(J6-Pin_4)"
$bom.Range("AC17").Value = "Net-(J6-Pin_1),Net-(J6-Pin_6),Net-(J6-Pin_2),Net-(J6-Pin_5),Net-(J6-Pin_3),Net-(J6-Pin_4)"

# Row 21 - big aggregated net list + its shorter label
$bom.Range("AB21").Value = "+5V,Net-(J3-Pin_4),Net-(J6-Pin_6),Net-(J4-Pin_1),unconnected-(U1-PB0-Pad14),Net-(J4-Pin_2),Net-(U1-PC0{slash}XTAL2),/TXLED,Earth,Net-(J6-Pin_4),Net-(U1-D+),Net-(J4-Pin_3),/MOSI2,Net-(U1-XTAL1),Net-(U1-UCAP),/RXLED,Net-(U1-D-),VBUS,Net-(J3-Pin_5),/MISO2,Net-(J3-Pin_1),Net-(J6-Pin_5),Net-(J4-Pin_4),Net-(J6-Pin_3),/DTR,Net-(J3-Pin_2),/RESET2,/SCK2,Net-(J6-Pin_2),Net-(J3-Pin_3),GND"
$bom.Range("AC21").Value = "SCK2,Net-(J6-Pin_2),Net-(J3-Pin_3),GND"

# Row 21 also grew taller (105pt -> 120pt)
$bom.Rows.Item(21).RowHeight = 120

# --- DNF sheet: update "Net Name" / "Net Label" columns (AB/AC) ---

# Row 9
$dnf.Range("AB9").Value = "Net-(U1-XTAL1),GND"
$dnf.Range("AC9").Value = "Net-(U1-XTAL1),GND"

# Row 11
$dnf.Range("AB11").Value = "Net-(J2-Shield),Earth"
$dnf.Range("AC11").Value = "Net-(J2-Shield),Earth"

# Row 12
$dnf.Range("AB12").Value = "Net-(J4-Pin_1),Net-(J4-Pin_2),Net-(J4-Pin_3),Net-(J4-Pin_4)"
$dnf.Range("AC12").Value = "Net-(J4-Pin_1),Net-(J4-Pin_2),Net-(J4-Pin_3),Net-(J4-Pin_4)"

# Row 13
$dnf.Range("AB13").Value = "Net-(J2-D-),Net-(J2-Shield),Earth,Net-(J2-D+),Net-(J2-VBUS)"
$dnf.Range("AC13").Value = "Net-(J2-D-),Net-(J2-Shield),Earth,Net-(J2-D+),Net-(J2-VBUS)"

# Row 15
$dnf.Range("AB15").Value = "Net-(U1-XTAL1),Net-(U1-PC0{slash}XTAL2)"
$dnf.Range("AC15").Value = "Net-(U1-XTAL1),Net-(U1-PC0{slash}XTAL2)"

# Row 16
$dnf.Range("AB16").Value = "Net-(J2-Shield),Net-(J2-D+)"
$dnf.Range("AC16").Value = "Net-(J2-Shield),Net-(J2-D+)"

# Row 17
$dnf.Range("AB17").Value = "Net-(U1-XTAL1),Net-(U1-PC0{slash}XTAL2)"
$dnf.Range("AC17").Value = "Net-(U1-XTAL1),Net-(U1-PC0{slash}XTAL2)"
